$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new (blank) column before column N ---
# This shifts the old column N ("Late") -> O, old O ("heading") -> P, and old
# P ("Outstanding") -> Q, leaving a new, empty column N in their place.
$ws = $wb.Worksheets.Item("Repayment schedule")
[void]$ws.Columns("N:N").Insert()

# Make "Repayment schedule" the active sheet/tab and select the cell the
# author ended up on after performing the column insert. Selecting here
# (rather than on "Transactions") is also what clears the previously
# active "tabSelected" flag on the "Transactions" sheet.
$ws.Activate()
[void]$ws.Range("P6").Select()
